$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '51.685.55'
$ws.Cells.Item(2, 5).Value = '  +1.37%  '

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.033.07'
$ws.Cells.Item(3, 5).Value = '  +2.68%  '

$ws.Cells.Item(4, 5).Value = '  +0.09%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '380.70'
$ws.Cells.Item(5, 5).Value = '  +0.55%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '102.93'
$ws.Cells.Item(6, 5).Value = '  +1.25%  '

$ws.Cells.Item(7, 5).Value = '  +0.75%  '

$ws.Cells.Item(9, 5).Value = '  +1.87%  '

$ws.Cells.Item(10, 5).Value = '  +1.57%  '

$ws.Cells.Item(11, 5).Value = '  -0.21%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.0860'
$ws.Cells.Item(12, 5).Value = '  +1.32%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '3.514.45'
$ws.Cells.Item(13, 5).Value = '  +2.77%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '18.57'
$ws.Cells.Item(14, 5).Value = '  +1.21%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.74'
$ws.Cells.Item(15, 5).Value = '  -0.69%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '3.038.27'
$ws.Cells.Item(16, 5).Value = '  +2.77%  '

$ws.Cells.Item(17, 5).Value = '  -3.37%  '

$ws.Cells.Item(18, 5).Value = '  -14.61%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '51.695.57'
$ws.Cells.Item(19, 5).Value = '  +1.42%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '3.06'
$ws.Cells.Item(20, 5).Value = '  -0.79%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '12.51'
$ws.Cells.Item(21, 5).Value = '  +0.98%  '

$ws.Cells.Item(22, 5).Value = '  +1.04%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '70.13'
$ws.Cells.Item(23, 5).Value = '  +0.98%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '268.50'
$ws.Cells.Item(24, 5).Value = '  +0.76%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '3.16'
$ws.Cells.Item(25, 5).Value = '  -1.40%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '8.28'
$ws.Cells.Item(26, 5).Value = '  +2.04%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '7.60'
$ws.Cells.Item(27, 5).Value = '  +8.01%  '

$ws.Cells.Item(28, 5).Value = '  +6.08%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 2).Value = 'Dai'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(29, 4).Value = '0.999'
$ws.Cells.Item(29, 5).Value = '  -0.10%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 2).Value = 'EthereumClassic'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(30, 4).Value = '26.26'
$ws.Cells.Item(30, 5).Value = '  +2.35%  '

$ws.Cells.Item(31, 5).Value = '  +1.04%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '10.28'
$ws.Cells.Item(32, 5).Value = '  +0.82%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '34.07'
$ws.Cells.Item(33, 5).Value = '  +0.91%  '

$ws.Cells.Item(34, 5).Value = '  +0.05%  '

$ws.Cells.Item(35, 5).Value = '  -0.15%  '

$ws.Cells.Item(36, 5).Value = '  +3.26%  '

$ws.Cells.Item(37, 5).Value = '  -0.08%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '3.34'
$ws.Cells.Item(38, 5).Value = '  +7.13%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.293'
$ws.Cells.Item(39, 5).Value = '  +13.52%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '17.08'
$ws.Cells.Item(40, 5).Value = '  +3.04%  '

$ws.Cells.Item(41, 5).Value = '  +2.56%  '

$ws.Cells.Item(42, 5).Value = '  +2.02%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 2).Value = 'Stellar'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(43, 4).Value = '0.116'
$ws.Cells.Item(43, 5).Value = '  -0.34%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 2).Value = 'NEARProtocol'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(44, 4).Value = '3.76'
$ws.Cells.Item(44, 5).Value = '  +6.29%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 2).Value = 'Monero'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(45, 4).Value = '123.99'
$ws.Cells.Item(45, 5).Value = '  +4.98%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '21.87'
$ws.Cells.Item(46, 5).Value = '  +2.40%  '

$ws.Cells.Item(47, 5).Value = '  +3.90%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.40'
$ws.Cells.Item(48, 5).Value = '  +3.84%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '2.035.53'
$ws.Cells.Item(49, 5).Value = '  +1.47%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '3.333.98'
$ws.Cells.Item(50, 5).Value = '  +2.84%  '

$ws.Cells.Item(51, 5).Value = '  +0.59%  '
